$d = $word.ActiveDocument

# Ordered list of (old, new) replacements exactly as they appear, in document
# order. Using wdReplaceOne (1) from the start of the document on each
# iteration guarantees that, even when the same "old" text occurs more than
# once (e.g. "25÷5="), each occurrence is replaced independently with the
# correct corresponding "new" value, matching document order.
$pairs = @(
    @("36÷4=", "97÷9="),
    @("38÷2=", "22÷5="),
    @("12÷2=", "92÷2="),
    @("91÷7=", "21÷9="),
    @("90÷9=", "51÷4="),
    @("64÷4=", "71÷7="),
    @("62÷4=", "35÷7="),
    @("80÷3=", "99÷8="),
    @("77÷8=", "12÷8="),
    @("34÷7=", "88÷9="),
    @("25÷5=", "30÷2="),
    @("58÷3=", "22÷9="),
    @("68÷3=", "11÷9="),
    @("79÷2=", "34÷8="),
    @("89÷2=", "86÷2="),
    @("10÷8=", "28÷3="),
    @("24÷4=", "87÷8="),
    @("24÷9=", "64÷3="),
    @("28÷9=", "24÷2="),
    @("48÷5=", "72÷6="),
    @("14÷5=", "47÷4="),
    @("25÷5=", "91÷2="),
    @("16÷8=", "35÷7="),
    @("28÷4=", "37÷8="),
    @("73÷8=", "44÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]

    $rng = $d.Content
    $rng.Start = 0

    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $new, 1)

    if (-not $found) {
        Write-Host "WARNING: could not find text to replace:" $old
    }
}
